# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.427.75'
$ws.Range("E2").Value = '  +1.91%  '

# Row 3
$ws.Range("D3").Value = '2.237.57'
$ws.Range("E3").Value = '  +0.93%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").Value = '''317.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.20%  '

# Row 6
$ws.Range("D6").Value = '''99.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.98%  '

# Row 7
$ws.Range("E7").Value = '  +2.02%  '

# Row 9
$ws.Range("D9").Value = '''0.563'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.16%  '

# Row 10
$ws.Range("D10").Value = '''37.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.19%  '

# Row 11
$ws.Range("E11").Value = '  -0.09%  '

# Row 12
$ws.Range("D12").Value = '''7.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.76%  '

# Row 13
$ws.Range("E13").Value = '  +2.55%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.579.87'
$ws.Range("E14").Value = '  +0.89%  '

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.867'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.30%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''14.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.12%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.260.38'
$ws.Range("E17").Value = '  +2.69%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '43.354.91'
$ws.Range("E18").Value = '  +1.95%  '

# Row 19
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").Value = '''14.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.93%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0975'
$ws.Range("E20").Value = '  +3.35%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''6.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.87%  '

# Row 22
$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").Value = '''65.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.41%  '

# Row 23
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '''3.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.22%  '

# Row 24
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '''236.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.10%  '

# Row 25
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = '''2.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.87%  '

# Row 26
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.08%  '

# Row 27
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '''4.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.73%  '

# Row 28
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '''10.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.55%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.29%  '

# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''6.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.88%  '

# Row 31
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '''36.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.76%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '''20.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.39%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.0872'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.28%  '

# Row 34
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Value = '''158.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.73%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''2.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.26%  '

# Row 36
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '''3.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.76%  '

# Row 37
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").Value = '''0.121'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.07%  '

# Row 38
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '''1.90'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.03%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''4.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.18%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.104'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '

# Row 41
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").Value = '''3.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.38%  '

# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0322'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.74%  '

# Row 43
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '''14.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +22.25%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.14%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.835.82'
$ws.Range("E45").Value = '  +3.56%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '''0.204'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.27%  '

# Row 47
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").Value = '''84.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.58%  '

# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '''5.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.21%  '

# Row 49
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''8.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.66%  '

# Row 50
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = '''74.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.20%  '

# Row 51
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '''58.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.42%  '

